$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the F3 cell (DBQ_Logical_Combo value for the "Five centimeters or
# larger in diameter" row) per the PR change - removes both the value and
# the style applied to that cell.
$ws.Range("F3").Clear()

# Leave the selection where the editor last clicked before saving.
$ws.Range("E5").Select() | Out-Null
